# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Add a new "metadata" worksheet, positioned right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"
$meta.Outline.SummaryRow = 1
$meta.Outline.SummaryColumn = 1

# Header row (B1:G1) — reuse the bold/bordered header style already used
# by the "data" sheet's own header row instead of building a new one.
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2) — A2 mirrors the bold index-column style from "data"
$meta.Cells.Item(2, 1).Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Cells.Item(2, 2).Value = "Primary ciliary disorders"
$meta.Cells.Item(2, 3).Value = 178
$meta.Cells.Item(2, 4).Value = "'1.32"
$meta.Cells.Item(2, 4).Style = "Normal"
$meta.Cells.Item(2, 5).Value = "2021-08-25T11:19:36.477951Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:22:17.515357"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/178/?format=json"

# --- Refresh the "time_taken" timestamps (column F) on the "data" sheet ---
$data.Cells.Item(2, 6).Value = "2021-10-05 14:22:17.517793"
$data.Cells.Item(3, 6).Value = "2021-10-05 14:22:17.517799"
$data.Cells.Item(4, 6).Value = "2021-10-05 14:22:17.517801"
$data.Cells.Item(5, 6).Value = "2021-10-05 14:22:17.517803"
$data.Cells.Item(6, 6).Value = "2021-10-05 14:22:17.517805"
$data.Cells.Item(7, 6).Value = "2021-10-05 14:22:17.517807"
$data.Cells.Item(8, 6).Value = "2021-10-05 14:22:17.517809"
$data.Cells.Item(9, 6).Value = "2021-10-05 14:22:17.517811"
$data.Cells.Item(10, 6).Value = "2021-10-05 14:22:17.517813"
$data.Cells.Item(11, 6).Value = "2021-10-05 14:22:17.517815"
$data.Cells.Item(12, 6).Value = "2021-10-05 14:22:17.517817"
$data.Cells.Item(13, 6).Value = "2021-10-05 14:22:17.517818"
$data.Cells.Item(14, 6).Value = "2021-10-05 14:22:17.517820"
$data.Cells.Item(15, 6).Value = "2021-10-05 14:22:17.517822"
$data.Cells.Item(16, 6).Value = "2021-10-05 14:22:17.517824"
$data.Cells.Item(17, 6).Value = "2021-10-05 14:22:17.517826"
$data.Cells.Item(18, 6).Value = "2021-10-05 14:22:17.517828"
$data.Cells.Item(19, 6).Value = "2021-10-05 14:22:17.517830"
$data.Cells.Item(20, 6).Value = "2021-10-05 14:22:17.517832"
$data.Cells.Item(21, 6).Value = "2021-10-05 14:22:17.517834"
$data.Cells.Item(22, 6).Value = "2021-10-05 14:22:17.517836"
$data.Cells.Item(23, 6).Value = "2021-10-05 14:22:17.517838"
$data.Cells.Item(24, 6).Value = "2021-10-05 14:22:17.517840"
$data.Cells.Item(25, 6).Value = "2021-10-05 14:22:17.517842"
$data.Cells.Item(26, 6).Value = "2021-10-05 14:22:17.517844"
$data.Cells.Item(27, 6).Value = "2021-10-05 14:22:17.517846"
$data.Cells.Item(28, 6).Value = "2021-10-05 14:22:17.517848"
$data.Cells.Item(29, 6).Value = "2021-10-05 14:22:17.517850"
$data.Cells.Item(30, 6).Value = "2021-10-05 14:22:17.517852"
$data.Cells.Item(31, 6).Value = "2021-10-05 14:22:17.517854"
$data.Cells.Item(32, 6).Value = "2021-10-05 14:22:17.517855"
$data.Cells.Item(33, 6).Value = "2021-10-05 14:22:17.517857"
$data.Cells.Item(34, 6).Value = "2021-10-05 14:22:17.517860"
$data.Cells.Item(35, 6).Value = "2021-10-05 14:22:17.517862"
$data.Cells.Item(36, 6).Value = "2021-10-05 14:22:17.517864"
$data.Cells.Item(37, 6).Value = "2021-10-05 14:22:17.517865"
$data.Cells.Item(38, 6).Value = "2021-10-05 14:22:17.517867"
$data.Cells.Item(39, 6).Value = "2021-10-05 14:22:17.517869"
$data.Cells.Item(40, 6).Value = "2021-10-05 14:22:17.517871"
$data.Cells.Item(41, 6).Value = "2021-10-05 14:22:17.517873"
$data.Cells.Item(42, 6).Value = "2021-10-05 14:22:17.517875"
$data.Cells.Item(43, 6).Value = "2021-10-05 14:22:17.517877"
$data.Cells.Item(44, 6).Value = "2021-10-05 14:22:17.517879"
$data.Cells.Item(45, 6).Value = "2021-10-05 14:22:17.517881"
$data.Cells.Item(46, 6).Value = "2021-10-05 14:22:17.517883"
$data.Cells.Item(47, 6).Value = "2021-10-05 14:22:17.517885"
$data.Cells.Item(48, 6).Value = "2021-10-05 14:22:17.517887"
$data.Cells.Item(49, 6).Value = "2021-10-05 14:22:17.517889"
$data.Cells.Item(50, 6).Value = "2021-10-05 14:22:17.517891"
$data.Cells.Item(51, 6).Value = "2021-10-05 14:22:17.517893"
$data.Cells.Item(52, 6).Value = "2021-10-05 14:22:17.517895"
$data.Cells.Item(53, 6).Value = "2021-10-05 14:22:17.517897"
$data.Cells.Item(54, 6).Value = "2021-10-05 14:22:17.517899"
$data.Cells.Item(55, 6).Value = "2021-10-05 14:22:17.517901"
$data.Cells.Item(56, 6).Value = "2021-10-05 14:22:17.517903"
$data.Cells.Item(57, 6).Value = "2021-10-05 14:22:17.517905"
$data.Cells.Item(58, 6).Value = "2021-10-05 14:22:17.517907"
$data.Cells.Item(59, 6).Value = "2021-10-05 14:22:17.517908"
$data.Cells.Item(60, 6).Value = "2021-10-05 14:22:17.517910"
$data.Cells.Item(61, 6).Value = "2021-10-05 14:22:17.517912"
$data.Cells.Item(62, 6).Value = "2021-10-05 14:22:17.517914"
$data.Cells.Item(63, 6).Value = "2021-10-05 14:22:17.517916"
$data.Cells.Item(64, 6).Value = "2021-10-05 14:22:17.517918"
$data.Cells.Item(65, 6).Value = "2021-10-05 14:22:17.517920"
$data.Cells.Item(66, 6).Value = "2021-10-05 14:22:17.517923"
$data.Cells.Item(67, 6).Value = "2021-10-05 14:22:17.517925"
$data.Cells.Item(68, 6).Value = "2021-10-05 14:22:17.517927"
$data.Cells.Item(69, 6).Value = "2021-10-05 14:22:17.517929"
$data.Cells.Item(70, 6).Value = "2021-10-05 14:22:17.517931"
$data.Cells.Item(71, 6).Value = "2021-10-05 14:22:17.517933"
$data.Cells.Item(72, 6).Value = "2021-10-05 14:22:17.517935"
$data.Cells.Item(73, 6).Value = "2021-10-05 14:22:17.517936"
$data.Cells.Item(74, 6).Value = "2021-10-05 14:22:17.517938"
$data.Cells.Item(75, 6).Value = "2021-10-05 14:22:17.517940"
$data.Cells.Item(76, 6).Value = "2021-10-05 14:22:17.517942"
$data.Cells.Item(77, 6).Value = "2021-10-05 14:22:17.517944"
$data.Cells.Item(78, 6).Value = "2021-10-05 14:22:17.517947"
$data.Cells.Item(79, 6).Value = "2021-10-05 14:22:17.517949"
$data.Cells.Item(80, 6).Value = "2021-10-05 14:22:17.517951"
$data.Cells.Item(81, 6).Value = "2021-10-05 14:22:17.517953"
$data.Cells.Item(82, 6).Value = "2021-10-05 14:22:17.517955"
$data.Cells.Item(83, 6).Value = "2021-10-05 14:22:17.517957"
$data.Cells.Item(84, 6).Value = "2021-10-05 14:22:17.517959"
$data.Cells.Item(85, 6).Value = "2021-10-05 14:22:17.517961"
$data.Cells.Item(86, 6).Value = "2021-10-05 14:22:17.517963"
$data.Cells.Item(87, 6).Value = "2021-10-05 14:22:17.517965"
$data.Cells.Item(88, 6).Value = "2021-10-05 14:22:17.517967"
$data.Cells.Item(89, 6).Value = "2021-10-05 14:22:17.517968"
$data.Cells.Item(90, 6).Value = "2021-10-05 14:22:17.517970"
$data.Cells.Item(91, 6).Value = "2021-10-05 14:22:17.517972"
$data.Cells.Item(92, 6).Value = "2021-10-05 14:22:17.517974"
$data.Cells.Item(93, 6).Value = "2021-10-05 14:22:17.517976"
$data.Cells.Item(94, 6).Value = "2021-10-05 14:22:17.517979"
$data.Cells.Item(95, 6).Value = "2021-10-05 14:22:17.517981"
$data.Cells.Item(96, 6).Value = "2021-10-05 14:22:17.517983"
$data.Cells.Item(97, 6).Value = "2021-10-05 14:22:17.517985"
$data.Cells.Item(98, 6).Value = "2021-10-05 14:22:17.517987"
$data.Cells.Item(99, 6).Value = "2021-10-05 14:22:17.517989"
$data.Cells.Item(100, 6).Value = "2021-10-05 14:22:17.517991"
$data.Cells.Item(101, 6).Value = "2021-10-05 14:22:17.517993"
$data.Cells.Item(102, 6).Value = "2021-10-05 14:22:17.517995"
$data.Cells.Item(103, 6).Value = "2021-10-05 14:22:17.517996"
$data.Cells.Item(104, 6).Value = "2021-10-05 14:22:17.517998"
$data.Cells.Item(105, 6).Value = "2021-10-05 14:22:17.518000"
$data.Cells.Item(106, 6).Value = "2021-10-05 14:22:17.518002"
$data.Cells.Item(107, 6).Value = "2021-10-05 14:22:17.518004"
$data.Cells.Item(108, 6).Value = "2021-10-05 14:22:17.518006"
$data.Cells.Item(109, 6).Value = "2021-10-05 14:22:17.518008"
$data.Cells.Item(110, 6).Value = "2021-10-05 14:22:17.518011"
$data.Cells.Item(111, 6).Value = "2021-10-05 14:22:17.518014"
$data.Cells.Item(112, 6).Value = "2021-10-05 14:22:17.518016"
$data.Cells.Item(113, 6).Value = "2021-10-05 14:22:17.518018"
$data.Cells.Item(114, 6).Value = "2021-10-05 14:22:17.518020"
$data.Cells.Item(115, 6).Value = "2021-10-05 14:22:17.518022"
$data.Cells.Item(116, 6).Value = "2021-10-05 14:22:17.518024"
$data.Cells.Item(117, 6).Value = "2021-10-05 14:22:17.518026"
$data.Cells.Item(118, 6).Value = "2021-10-05 14:22:17.518028"
$data.Cells.Item(119, 6).Value = "2021-10-05 14:22:17.518030"
$data.Cells.Item(120, 6).Value = "2021-10-05 14:22:17.518032"
$data.Cells.Item(121, 6).Value = "2021-10-05 14:22:17.518034"
$data.Cells.Item(122, 6).Value = "2021-10-05 14:22:17.518036"
$data.Cells.Item(123, 6).Value = "2021-10-05 14:22:17.518038"
$data.Cells.Item(124, 6).Value = "2021-10-05 14:22:17.518040"
$data.Cells.Item(125, 6).Value = "2021-10-05 14:22:17.518042"
$data.Cells.Item(126, 6).Value = "2021-10-05 14:22:17.518044"
$data.Cells.Item(127, 6).Value = "2021-10-05 14:22:17.518046"
$data.Cells.Item(128, 6).Value = "2021-10-05 14:22:17.518048"
$data.Cells.Item(129, 6).Value = "2021-10-05 14:22:17.518051"
$data.Cells.Item(130, 6).Value = "2021-10-05 14:22:17.518054"
$data.Cells.Item(131, 6).Value = "2021-10-05 14:22:17.518057"
$data.Cells.Item(132, 6).Value = "2021-10-05 14:22:17.518059"
$data.Cells.Item(133, 6).Value = "2021-10-05 14:22:17.518061"
$data.Cells.Item(134, 6).Value = "2021-10-05 14:22:17.518063"
$data.Cells.Item(135, 6).Value = "2021-10-05 14:22:17.518065"
$data.Cells.Item(136, 6).Value = "2021-10-05 14:22:17.518067"
$data.Cells.Item(137, 6).Value = "2021-10-05 14:22:17.518069"
$data.Cells.Item(138, 6).Value = "2021-10-05 14:22:17.518071"
$data.Cells.Item(139, 6).Value = "2021-10-05 14:22:17.518073"
$data.Cells.Item(140, 6).Value = "2021-10-05 14:22:17.518075"
$data.Cells.Item(141, 6).Value = "2021-10-05 14:22:17.518077"
$data.Cells.Item(142, 6).Value = "2021-10-05 14:22:17.518079"
$data.Cells.Item(143, 6).Value = "2021-10-05 14:22:17.518081"
$data.Cells.Item(144, 6).Value = "2021-10-05 14:22:17.518083"
